# The tram-17 station retrieval routine had dropped two real stops
# ("Jan Pieter Heijestraat" and "Ten Katestraat") and a stray, incorrect
# one ("Nicolaas Beetsstraat") had been recorded in their place. Fix the
# list in column A: remove the wrong row and insert the two correct
# stations between "Witte de Withstraat" (row 13) and "Bilderdijkstraat"
# (row 14 -> becomes row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 14 ("Nicolaas Beetsstraat") and everything below it down by
# one row, opening up two blank rows at 14 and 15.
$ws.Rows.Item(14).Insert()

# Fill the two newly opened rows with the correct station names. The old
# "Nicolaas Beetsstraat" value (now shifted to row 15) is overwritten and
# disappears from the workbook entirely.
$ws.Range("A14").Value = "Jan Pieter Heijestraat"
$ws.Range("A15").Value = "Ten Katestraat"

# Leave the selection where the author's saved file shows it.
$ws.Range("D17").Select() | Out-Null
